$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear all existing hyperlinks on the sheet (will re-add for the final row set)
$ws.Range("F2").Hyperlinks.Delete()

# Clear the previously-used data rows (2-10) plus leftover columns so no stale values remain
$ws.Range("A2:H10").ClearContents()

# Row 2
$ws.Cells.Item(2,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(2,2).Value = "【急募】LINEとChatGPT連携の簡易質問対応システム開発"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5452848"
$ws.Cells.Item(2,7).Value = 430
$ws.Cells.Item(2,8).Value = "🔥GPT,ChatGPT ◆開発,システム開発"
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://www.lancers.jp/work/detail/5452848")

# Row 3
$ws.Cells.Item(3,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(3,2).Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Cells.Item(3,7).Value = 385
$ws.Cells.Item(3,8).Value = "🔥AI,Ai ◆効率化"
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), "https://www.lancers.jp/work/detail/5423720")

# Row 4
$ws.Cells.Item(4,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(4,2).Value = "初回 スポーツクラブ コスパ自動予約bot開発(playwight/Python)"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5452614"
$ws.Cells.Item(4,7).Value = 368
$ws.Cells.Item(4,8).Value = "🔥Python ★bot ◆開発"
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), "https://www.lancers.jp/work/detail/5452614")

# Row 5
$ws.Cells.Item(5,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(5,2).Value = "AIオートメーションエンジニア"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5452520"
$ws.Cells.Item(5,7).Value = 303
$ws.Cells.Item(5,8).Value = "🔥AI,Ai"
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), "https://www.lancers.jp/work/detail/5452520")

# Row 6
$ws.Cells.Item(6,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(6,2).Value = "【Flutterエンジニア募集】Androidアプリ開発のパートナーを探しています"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5452211"
$ws.Cells.Item(6,7).Value = 100
$ws.Cells.Item(6,8).Value = "◆開発 ◇アプリ"
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), "https://www.lancers.jp/work/detail/5452211")

# Row 7
$ws.Cells.Item(7,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(7,2).Value = "製造業向けMR業務支援アプリケーションの機能開発エンジニア募集(Unity/C#)"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5441557"
$ws.Cells.Item(7,7).Value = 93
$ws.Cells.Item(7,8).Value = "◆開発 ◇アプリ"
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5441557")

# Row 8
$ws.Cells.Item(8,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(8,2).Value = "【再掲】基幹システム入替に伴うBIツール環境の再構築(Microsoft Power BI)"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5452367"
$ws.Cells.Item(8,7).Value = 88
$ws.Cells.Item(8,8).Value = "◆ツール"
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), "https://www.lancers.jp/work/detail/5452367")

# Row 9
$ws.Cells.Item(9,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(9,2).Value = "製造業向け 技能習得・作業トレーニングVRシステムの開発(Unity/R3)"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5441568"
$ws.Cells.Item(9,7).Value = 83
$ws.Cells.Item(9,8).Value = "◆開発"
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), "https://www.lancers.jp/work/detail/5441568")

# Row 10
$ws.Cells.Item(10,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(10,2).Value = "Amazonの購入アカウントから必要な情報のスクレイピング→スプレッドシートに記入をしたい。"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5452210"
$ws.Cells.Item(10,7).Value = 40
$ws.Cells.Item(10,8).Value = "◆スクレイピング"
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), "https://www.lancers.jp/work/detail/5452210")

# Row 11
$ws.Cells.Item(11,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(11,2).Value = "【R/Shiny】高齢者評価アプリ 機能追加・UI改修依頼"
$ws.Cells.Item(11,3).Value = "システム開発"
$ws.Cells.Item(11,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(11,5).Value = "期限情報なし"
$ws.Cells.Item(11,6).Value = "https://www.lancers.jp/work/detail/5452159"
$ws.Cells.Item(11,7).Value = 38
$ws.Cells.Item(11,8).Value = "◇アプリ"
$ws.Hyperlinks.Add($ws.Cells.Item(11,6), "https://www.lancers.jp/work/detail/5452159")

# Row 12
$ws.Cells.Item(12,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(12,2).Value = "【小規模・短納期・急募】アプリMatrixifyを用いたデータ移行検証・マッピング担当募集"
$ws.Cells.Item(12,3).Value = "システム開発"
$ws.Cells.Item(12,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(12,5).Value = "期限情報なし"
$ws.Cells.Item(12,6).Value = "https://www.lancers.jp/work/detail/5451926"
$ws.Cells.Item(12,7).Value = 33
$ws.Cells.Item(12,8).Value = "◇アプリ"
$ws.Hyperlinks.Add($ws.Cells.Item(12,6), "https://www.lancers.jp/work/detail/5451926")

# Row 13
$ws.Cells.Item(13,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(13,2).Value = "注目 限定公開 PR 限定公開の仕事"
$ws.Cells.Item(13,3).Value = "システム開発"
$ws.Cells.Item(13,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(13,5).Value = "期限情報なし"
$ws.Cells.Item(13,6).Value = "https://www.lancers.jp/work/detail/5450323"
$ws.Cells.Item(13,7).Value = 13
$ws.Hyperlinks.Add($ws.Cells.Item(13,6), "https://www.lancers.jp/work/detail/5450323")

# Row 14
$ws.Cells.Item(14,1).Value = "2025-12-12 18:29:17"
$ws.Cells.Item(14,2).Value = "Xの運用代行"
$ws.Cells.Item(14,3).Value = "システム開発"
$ws.Cells.Item(14,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(14,5).Value = "期限情報なし"
$ws.Cells.Item(14,6).Value = "https://www.lancers.jp/work/detail/5451931"
$ws.Cells.Item(14,7).Value = 10
$ws.Hyperlinks.Add($ws.Cells.Item(14,6), "https://www.lancers.jp/work/detail/5451931")

# Widen column H (skill summary) to fit the longer new entries
$ws.Columns.Item(8).ColumnWidth = 24.15

